$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> (new Price text, new Volume(1h) text). $null means "leave unchanged".
$updates = @{
    2 = @('27.912.67', '  -5.56%  ')
    3 = @('1.819.78', '  -4.59%  ')
    4 = @('1.004', '  -0.65%  ')
    5 = @('329.16', '  -2.68%  ')
    6 = @('1.003', '  -0.61%  ')
    7 = @('0.4629', '  -2.76%  ')
    8 = @('0.3841', '  -4.06%  ')
    9 = @('45.87', '  -3.70%  ')
    10 = @('0.07825', '  -2.69%  ')
    11 = @('0.9585', '  -3.08%  ')
    12 = @('21.86', '  -6.17%  ')
    13 = @('1.807.59', '  -5.60%  ')
    14 = @('5.637', '  -4.83%  ')
    15 = @('6.841', '  -3.90%  ')
    16 = @('0.06855', '  +0.32%  ')
    17 = @('1.002', '  -0.89%  ')
    18 = @('86.59', '  -2.89%  ')
    19 = @('0.000009938', '  -2.55%  ')
    20 = @('16.68', '  -4.03%  ')
    21 = @('1.003', '  -0.82%  ')
    22 = @('27.957.19', '  -5.47%  ')
    23 = @('5.311', '  -3.65%  ')
    24 = @('10.91', '  -6.12%  ')
    25 = @('2.101', '  -2.43%  ')
    26 = @('2.074.51', '  -3.25%  ')
    27 = @('152.72', '  -2.66%  ')
    28 = @('19.16', '  -2.57%  ')
    29 = @('5.705', '  -12.58%  ')
    30 = @('1.960', '  -4.63%  ')
    31 = @($null, '  -2.42%  ')
    32 = @('0.9367', '  -5.75%  ')
    33 = @('0.09250', '  -2.85%  ')
    34 = @('5.265', '  -4.10%  ')
    35 = @('3.422', '  -3.63%  ')
    36 = @($null, '  -5.70%  ')
    37 = @('0.05942', '  -7.98%  ')
    38 = @($null, '  -4.39%  ')
    39 = @('1.148', '  -3.80%  ')
    40 = @($null, '  -0.87%  ')
    41 = @('7.564', '  -2.60%  ')
    42 = @('0.5578', '  -4.25%  ')
    43 = @('9.887', '  -6.19%  ')
    44 = @('0.1768', '  -2.91%  ')
    45 = @('1.219', '  -4.23%  ')
    46 = @('2.219', '  -9.18%  ')
    47 = @('11.56', '  -5.19%  ')
    48 = @('0.5244', '  -4.19%  ')
    49 = @('0.07001', '  -5.77%  ')
    50 = @('1.821', '  -6.11%  ')
    51 = @('112.50', '  -3.18%  ')
}

foreach ($row in $updates.Keys) {
    $priceText  = $updates[$row][0]
    $volumeText = $updates[$row][1]

    if ($priceText -ne $null) {
        # Column D ("Price") holds values like "27.912.67" or "0.07825" that Excel
        # would otherwise auto-parse as a number (and mangle via float rounding, or
        # reject outright for the thousand-dot style). Force the cell to Text first
        # so the literal string is preserved exactly, matching the inline-string cells
        # already used throughout this column.
        $priceCell = $ws.Cells.Item($row, 4)
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $priceText
    }

    if ($volumeText -ne $null) {
        # Column E ("Volume(1h)") values are padded percentages (e.g. "  -5.56%  ")
        # and already round-trip as text without help.
        $ws.Cells.Item($row, 5).Value = $volumeText
    }
}
